$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the previously entered item data (columns B:K) for the existing rows,
# keeping only the running "Sl.no." counter in column A. This also drops the
# now-unused shared strings that described those items.
$ws.Range("B2:K6").ClearContents()

# Add a new row (7) continuing the Sl.no. sequence, fixing the
# "add item" issue from the commit message.
$ws.Range("A7").Value = 6
